# Update the seating arrangement: reshuffle the "Occupant" column (C) for the
# existing tables (rows 2-25) and append two new tables (7 and 8) occupying
# rows 26-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows: only the Occupant name (column C) changes ---
$ws.Range("C2").Value = "Younes"
$ws.Range("C3").Value = "Nadiya"
$ws.Range("C4").Value = "Choti"
$ws.Range("C5").Value = "Sofia"
$ws.Range("C6").Value = "Jordi"
$ws.Range("C7").Value = "Marc5"
$ws.Range("C8").Value = "Hajer"
$ws.Range("C9").Value = "Charly"
$ws.Range("C10").Value = "Megan"
$ws.Range("C11").Value = "Kenny"
$ws.Range("C12").Value = "Dragos"
$ws.Range("C13").Value = "Fang"
$ws.Range("C14").Value = "Aida"
$ws.Range("C15").Value = "Marc2"
$ws.Range("C16").Value = "Emmanuel"
$ws.Range("C17").Value = "Augustin"
$ws.Range("C18").Value = "Mengstu"
$ws.Range("C19").Value = "Floriane"
$ws.Range("C20").Value = "Yassine"
$ws.Range("C21").Value = "Caterina"
$ws.Range("C22").Value = "Marc"
$ws.Range("C23").Value = "Marc3"
$ws.Range("C24").Value = "Yves"
$ws.Range("C25").Value = "Klebert"

# --- New rows: table 7 (rows 26-29) and table 8 (rows 30-33) ---
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "Hanieh"

$ws.Range("A27").Value = 7
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "Alberto"

$ws.Range("A28").Value = 7
$ws.Range("B28").Value = 3
$ws.Range("C28").Value = "Evi"

$ws.Range("A29").Value = 7
$ws.Range("B29").Value = 4
$ws.Range("C29").Value = "Moussa"

$ws.Range("A30").Value = 8
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "Marc4"

$ws.Range("A31").Value = 8
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = "Preeti"

$ws.Range("A32").Value = 8
$ws.Range("B32").Value = 3
$ws.Range("C32").Value = "Empty"

$ws.Range("A33").Value = 8
$ws.Range("B33").Value = 4
$ws.Range("C33").Value = "Empty"
